$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    3  = 7378000000.0
    4  = -1496000000.0
    5  = -31000000.0
    6  = 27566000000.0
    7  = 28642000000.0
    8  = -1282000000.0
    9  = 24252000000.0
    10 = -2508000000.0
    12 = -9193000000.0
    13 = -2559000000.0
    14 = -20429000000.0
    15 = 5662000000.0
    16 = -1930000000.0
    17 = -10635000000.0
    18 = 702000000.0
    19 = -6930000000.0
    20 = 248000000.0
    21 = -2859000000.0
    22 = 13985000000.0
    23 = 11126000000.0
    24 = 1049000000.0
    25 = -10635000000.0
    26 = 3255000000.0
    27 = -9193000000.0
    28 = -1930000000.0
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 2).Value = $values[$row]
}
